$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two time-range trigger cells
$ws.Range("B20").Value = "21:15 - 21:19"
$ws.Range("B21").Value = "21:20 - 21:24"

# Update the active selection in the sheet view
$ws.Range("B24").Select()
